$wb = $excel.ActiveWorkbook

$wsNodes = $wb.Worksheets.Item("nodes")
$wsInteractions = $wb.Worksheets.Item("interactions")
$wsUnits = $wb.Worksheets.Item("units")
$wsParameters = $wb.Worksheets.Item("parameters")

# ---------------------------------------------------------------------------
# "nodes" sheet: rename abbreviated node/group labels to descriptive ones.
#   other -> north_side, side -> south_side
#   b1/b2/b3 -> bridge_1/bridge_2/bridge_3
#   h1 -> hill_1, h2 -> hill_2, h3 & h4 -> hill_3
# ---------------------------------------------------------------------------
$wsNodes.Range("E2:E9").Value = "north_side"
$wsNodes.Range("E10:E12").Value = "south_side"
$wsNodes.Range("E14").Value = "south_side"

$wsNodes.Range("E13").Value = "hill_1"
$wsNodes.Range("E15:E16").Value = "hill_2"
$wsNodes.Range("E17:E19").Value = "hill_3"

$wsNodes.Range("F7").Value = "bridge_1"
$wsNodes.Range("F8").Value = "bridge_2"
$wsNodes.Range("F9").Value = "bridge_3"
$wsNodes.Range("F10").Value = "bridge_1"
$wsNodes.Range("F11").Value = "bridge_2"
$wsNodes.Range("F12").Value = "bridge_3"

# ---------------------------------------------------------------------------
# "interactions" sheet: same renames for the from/to node-pair columns.
# ---------------------------------------------------------------------------
$wsInteractions.Range("A2").Value = "bridge_1"
$wsInteractions.Range("B2").Value = "bridge_1"
$wsInteractions.Range("A3").Value = "bridge_2"
$wsInteractions.Range("B3").Value = "bridge_2"
$wsInteractions.Range("A4").Value = "bridge_3"
$wsInteractions.Range("B4").Value = "bridge_3"
$wsInteractions.Range("A5").Value = "bridge_1,bridge_2"
$wsInteractions.Range("B5").Value = "bridge_1,bridge_2"
$wsInteractions.Range("A6").Value = "bridge_2,bridge_3"
$wsInteractions.Range("B6").Value = "bridge_2,bridge_3"

# Highlight the towers/other interaction row with the same fill already used
# for "towers" nodes elsewhere in the workbook (copy format only).
$wsNodes.Range("F13").Copy()
$wsInteractions.Range("A9").PasteSpecial(-4122)

# Widen the first two columns so the longer node names are readable.
$wsInteractions.Columns("A:B").ColumnWidth = 24.7109375

# ---------------------------------------------------------------------------
# Restore each sheet's own last selection (captured from the authored file),
# finishing on "nodes" so it stays the active/visible tab.
# ---------------------------------------------------------------------------
$wsInteractions.Range("C23").Select()
$wsParameters.Range("B36").Select()
$wsUnits.Range("I38").Select()
$wsNodes.Range("D27").Select()
